$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (columns reordered/renamed): lang_code, code, name, descr, is_active
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# Data rows 2..11 : lang_code, code, name, descr, is_active
$data = @(
    @("eng", "FRS", "Finger Print Scanner", "For scanning fingerprints", $true),
    @("fra", "FRS", "Scanner dempreintes digitales", "Scannez les empreintes digitales", $true),
    @("eng", "IRS", "Iris Scanner", "For scanning Iris", $true),
    @("fra", "IRS", "Scanner dIris", "Pour scanner liris", $true),
    @("eng", "CMR", "Camera", "For capturing photo", $true),
    @("fra", "CMR", "Caméra", "Pour capturer une photo", $true),
    @("eng", "SCN", "Document Scanner", "For scanning documents", $true),
    @("fra", "SCN", "Scanner de documents", "Pour numériser des documents", $true),
    @("eng", "PRT", "Printer", "For printing Documents", $true),
    @("fra", "PRT", "Imprimante", "Pour imprimer des documents", $true)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row = $row + 1
}

# Column A (lang_code) uses the same bold/bordered style as the header row
# for every data row (A2:A11), matching the style applied to A1.
$ws.Range("A1").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
